$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (1): 0..7 across A1:H1
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7

# Row 2
$ws.Range("A2").Value = "inf"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 25
$ws.Range("E2").Value = 40
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = 65
$ws.Range("H2").Value = 75

# Row 3
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "inf"
$ws.Range("C3").Value = 9
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 34
$ws.Range("F3").Value = 45
$ws.Range("G3").Value = 58
$ws.Range("H3").Value = 72

# Row 4
$ws.Range("A4").Value = 15
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = "inf"
$ws.Range("D4").Value = 11
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = 37
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 64

# Row 5
$ws.Range("A5").Value = 25
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 11
$ws.Range("D5").Value = "inf"
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = 26
$ws.Range("G5").Value = 38
$ws.Range("H5").Value = 54

# Row 6
$ws.Range("A6").Value = 40
$ws.Range("B6").Value = 34
$ws.Range("C6").Value = 25
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = "inf"
$ws.Range("F6").Value = 13
$ws.Range("G6").Value = 26
$ws.Range("H6").Value = 42

# Row 7
$ws.Range("A7").Value = 50
$ws.Range("B7").Value = 45
$ws.Range("C7").Value = 37
$ws.Range("D7").Value = 26
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = "inf"
$ws.Range("G7").Value = 14
$ws.Range("H7").Value = 30

# Row 8
$ws.Range("A8").Value = 65
$ws.Range("B8").Value = 58
$ws.Range("C8").Value = 50
$ws.Range("D8").Value = 38
$ws.Range("E8").Value = 26
$ws.Range("F8").Value = 14
$ws.Range("G8").Value = "inf"
$ws.Range("H8").Value = 17

# Row 9
$ws.Range("A9").Value = 75
$ws.Range("B9").Value = 72
$ws.Range("C9").Value = 64
$ws.Range("D9").Value = 54
$ws.Range("E9").Value = 42
$ws.Range("F9").Value = 30
$ws.Range("G9").Value = 17
$ws.Range("H9").Value = "inf"

# Update the active selection to match the saved view state (D13)
$ws.Range("D13").Select()
